$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each row.
# D values are entered with a leading apostrophe to force text (avoids Excel
# auto-converting numeric-looking strings like "1.006" into actual numbers,
# which would also truncate values like "7.640" -> 7.64). The Style reset
# keeps the cell formatting identical to the original (no visible quote
# prefix indicator / no added number format).

$ws.Range("D2").Value = "'26.747.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.70%  "

$ws.Range("D3").Value = "'1.697.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.84%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.38%  "

$ws.Range("D5").Value = "'218.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.42%  "

$ws.Range("D6").Value = "'0.5053"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -15.21%  "

$ws.Range("D7").Value = "'1.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").Value = "'0.2599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.33%  "

$ws.Range("D9").Value = "'21.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.46%  "

$ws.Range("D10").Value = "'0.06122"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.52%  "

$ws.Range("D11").Value = "'0.07322"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.75%  "

$ws.Range("D12").Value = "'1.727.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.03%  "

$ws.Range("D13").Value = "'4.422"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.99%  "

$ws.Range("D14").Value = "'1.929.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.81%  "

$ws.Range("D15").Value = "'0.5724"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.63%  "

$ws.Range("D16").Value = "'0.000008163"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -11.16%  "

$ws.Range("D17").Value = "'65.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -13.31%  "

$ws.Range("D18").Value = "'26.793.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.47%  "

$ws.Range("D19").Value = "'5.015"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.06%  "

$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").Value = "'10.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.10%  "

$ws.Range("D22").Value = "'184.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -12.28%  "

$ws.Range("D23").Value = "'6.217"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -9.13%  "

$ws.Range("D24").Value = "'1.007"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").Value = "'145.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.85%  "

$ws.Range("D26").Value = "'7.640"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.56%  "

$ws.Range("D27").Value = "'0.1142"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.50%  "

$ws.Range("D28").Value = "'15.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.73%  "

$ws.Range("D29").Value = "'1.319"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.78%  "

$ws.Range("D30").Value = "'0.05628"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.80%  "

$ws.Range("D31").Value = "'1.331"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.18%  "

$ws.Range("D32").Value = "'3.467"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.85%  "

$ws.Range("D33").Value = "'3.439"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.73%  "

$ws.Range("D34").Value = "'1.659"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.93%  "

$ws.Range("D35").Value = "'1.005"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.28%  "

$ws.Range("E36").Value = "  -3.87%  "

$ws.Range("D37").Value = "'0.5895"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.22%  "

$ws.Range("D38").Value = "'2.629"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.21%  "

$ws.Range("D39").Value = "'0.01590"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.78%  "

$ws.Range("D40").Value = "'1.068.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.97%  "

$ws.Range("D41").Value = "'5.883"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.85%  "

$ws.Range("D42").Value = "'0.8517"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.38%  "

$ws.Range("D43").Value = "'1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").Value = "'98.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.18%  "

$ws.Range("D45").Value = "'1.859.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.11%  "

$ws.Range("D46").Value = "'56.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.68%  "

$ws.Range("D49").Value = "'8.102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.59%  "

$ws.Range("D50").Value = "'0.4337"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.62%  "

$ws.Range("D51").Value = "'0.05207"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.25%  "

# Rows 47 and 48 swap (Frax <-> BabyDogeCoin) along with their new price/volume data.
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "'1.008"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.00000000104"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.95%  "
